# "harmonized similar tags to be the same"
#
# The "isa_template" sheet has a Tags row (row 12) annotated with a term
# (row 12 = Tags value / term source REF columns, row 13 = Tags Term
# Accession Number / Tags Term Source REF). Previously the "Tags" row had
# three loosely related, unpaired tag values ("Measurement", "Mass
# spectrometry", "MS") with no accession numbers. This change harmonizes
# them into two properly paired Tag + Term-Accession-Number + Term-Source
# pairs: "measurement" (EFO:0001444) and "Mass Spectrometry"
# (NCIT:C17156).

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("isa_template")

# Row 12 ("Tags") / Row 13 ("Tags Term Accession Number"): rename/re-case
# the two tag terms and give each its proper ontology accession number,
# going column by column (Tags, then its Term Accession Number, for each
# of the two harmonized terms in turn).
$ws2.Range("C12").Value2 = "measurement"
$ws2.Range("C13").Value2 = "EFO:0001444"
$ws2.Range("D12").Value2 = "Mass Spectrometry"
$ws2.Range("D13").Value2 = "NCIT:C17156"

# Drop the stray 3rd tag ("MS") that had no matching accession number.
$ws2.Range("E12").Value2 = ""

# Row 13 now wraps onto a second line like row 12 above it, so its row
# grows to match.
$ws2.Rows("13").RowHeight = $ws2.Rows("12").RowHeight

# Reflect where the edit left the selection.
$ws2.Activate() | Out-Null
$ws2.Range("E20").Select() | Out-Null
